$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) "NOME COMPLETO DO ALUNO" cell: append the remaining group members
#    after "Valdiney Atílio Pedro".
# ------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("Valdiney Atílio Pedro", $true, $false, $false, $false, $false, `
                   $true, 1, $false, "", 0) | Out-Null
$rng.Collapse(0)
$rng.InsertAfter("; Patrícia Corrêa França; Mariana Alcantara; Mariana Simões.")

# ------------------------------------------------------------------
# 2) "RA" cell: append the remaining group RAs after "10424616".
# ------------------------------------------------------------------
$rng2 = $d.Content
$rng2.Find.Execute("10424616", $true, $false, $false, $false, $false, `
                    $true, 1, $false, "", 0) | Out-Null
$rng2.Collapse(0)
$rng2.InsertAfter("; 10423533; ; 10424388.")

# ------------------------------------------------------------------
# 3) Fix the mangled GitHub hyperlink display text (it had been split
#    across several runs with a stray bookmark in between).
# ------------------------------------------------------------------
$h = $d.Hyperlinks.Item(1)
$h.TextToDisplay = "https://github.com/valdineyatilio/Projeto-Aplicado-I"
